$d = $word.ActiveDocument

# 1. Remove the trailing space after "TradeTracker" and remove the word "app" (keep following space)
$d.Content.Find.Execute("TradeTracker app ", $true, $false, $false, $false, $false, $true, 1, $false, "TradeTracker  ", 2)

# 2. Add a trailing space after the sentence-ending period, and rewrite "(more details in projects)"
#    to "(details in projects section)"
$d.Content.Find.Execute("for android. (more details in projects)", $true, $false, $false, $false, $false, $true, 1, $false, "for android. (details in projects section)", 2)

# 3. Remove the old bookmark (it will be re-created at the new location)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 4. Insert the _GoBack bookmark at its new location: right after "details in" and before " projects section)"
$r = $d.Content
$r.Find.Execute("details in", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $anchor)
